$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 7: date entry + comment, matching the style of the existing
# date/comment rows (row 6: A6 uses numFmt 14 "mm-dd-yy" + vertical-top;
# D6 uses wrap-text + vertical-top).
$ws.Range("A7").Value = 42986
$ws.Range("A7").NumberFormat = "mm-dd-yy"
$ws.Range("A7").VerticalAlignment = -4160

$ws.Range("D7").Value = "[Objekt aufnehmen/werfen noch mit Bugs]; Menüs eingeführt für Pause/Gewonnen/Verloren"
$ws.Range("D7").WrapText = $true
$ws.Range("D7").VerticalAlignment = -4160

$ws.Rows.Item(7).RowHeight = 45

# Move the active selection down to the next empty row, like Excel does
# after entering data in the row above.
$ws.Range("A8").Select()
